# Set an explicit (custom) height of 24 points on row 1 of every worksheet,
# matching the "Power_BusInfo" example update (both ScenarioA and ScenarioB).
$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Rows.Item(1).RowHeight = 24
}
